# Add a new "BOUNDARY" attack block (columns AS:AZ) mirroring the existing
# attack blocks (REV, REV_NO_EQUAL, REV_BIM, FGSM_SURRO, FGSM) for SEED 888.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: block header "BOUNDARY" in AS1, merged across AS1:AZ1 -----------
# Merge first, then stamp the style from an existing header block (AK1) over
# the whole merged range, then set the visible text - this order keeps every
# cell in the merge on the same plain header style (s=1), matching how the
# pre-existing header blocks (e.g. AK1:AR1) look.
$ws.Range("AS1:AZ1").Merge()
$ws.Range("AK1").Copy()
$ws.Range("AS1:AZ1").PasteSpecial(-4122)
$ws.Range("AS1").Value = "BOUNDARY"

# --- Row 2: epsilon labels (stored as text, matching existing ε row) -------
$epsCols = @("AS","AT","AU","AV","AW","AX","AY","AZ")
$epsVals = @("0.01","0.02","0.03","0.04","0.05","0.07","0.10","0.20")
for ($i = 0; $i -lt $epsCols.Length; $i++) {
    $ws.Range($epsCols[$i] + "2").Value = "'" + $epsVals[$i]
}
$ws.Range("M2").Copy()
$ws.Range("AS2:AZ2").PasteSpecial(-4122)

# --- Rows 4-12: numeric metric values (no explicit style, like C:AR) -------
$dataCols = @("AS","AT","AU","AV","AW","AX","AY","AZ")

$row4 = @(4.225946426391602, 4.250112533569336, 4.336052894592285, 4.419761180877686, 4.485738754272461, 4.650653839111328, 4.967658996582031, 6.78828239440918)
$row5 = @(5.380455087749915, 5.393753691812359, 5.469000935928117, 5.600021512126475, 5.667467995626028, 5.902176678537989, 6.213391626491155, 8.565021227238546)
$row6 = @(0.9996083378791809, 0.9996046423912048, 0.9995924234390259, 0.9995657205581665, 0.9995605945587158, 0.9995237588882446, 0.9994716048240662, 0.9989577531814575)
$row7 = @(2.530349493026733, 2.62942099571228, 3.02424693107605, 3.525465488433838, 3.796053171157837, 5.369766712188721, 6.114477157592773, 11.37678623199463)
$row8 = @(3.397235385900117, 3.526560547523482, 3.989035360421044, 4.516633551406989, 4.837836809286371, 6.799186927206351, 7.795911343332202, 14.370735929223)
$row9 = @(0.9997932910919189, 0.9997768998146057, 0.9997138977050781, 0.999633252620697, 0.9995796084403992, 0.9991684556007385, 0.998913049697876, 0.9962786436080933)
$row10 = @(2.684278964996338, 2.818708896636963, 3.036305904388428, 3.248105525970459, 3.627833604812622, 4.270208358764648, 5.393082141876221, 9.396393775939941)
$row11 = @(3.573498128986309, 3.72261359559607, 3.946034220033302, 4.187992123965315, 4.542967700384122, 5.41309853431355, 6.670052558657534, 11.65110423867492)
$row12 = @(0.9997835159301758, 0.9997653365135193, 0.9997282028198242, 0.9996936321258545, 0.9996384978294373, 0.9994862079620361, 0.9992225766181946, 0.9976072311401367)

$allRows = @{ 4 = $row4; 5 = $row5; 6 = $row6; 7 = $row7; 8 = $row8; 9 = $row9; 10 = $row10; 11 = $row11; 12 = $row12 }

foreach ($r in @(4, 5, 6, 7, 8, 9, 10, 11, 12)) {
    $vals = $allRows[$r]
    for ($i = 0; $i -lt $dataCols.Length; $i++) {
        $ws.Range($dataCols[$i] + $r).Value = $vals[$i]
    }
}
